# "write some new for stephen"
#
# The original paragraph was split across three separate <w:r> runs that
# all share identical run formatting (sz=24, szCs=24, rtl=0):
#
#   "I ended up being late for class because of my trip, but thankfully instead "
#   "of a scolding"
#   " Ms. Tran simply told me to take better care of myself and continued on with the lesson."
#
# Clean this up into a single run holding the full, concatenated sentence.
# Because all three runs carry the same formatting, a plain Find/Replace
# across the whole sentence naturally collapses them into one run.

$d = $word.ActiveDocument

$oldText = "I ended up being late for class because of my trip, but thankfully instead of a scolding Ms. Tran simply told me to take better care of myself and continued on with the lesson."
$newText = "I ended up being late for class because of my trip, but thankfully instead of a scolding Ms. Tran simply told me to take better care of myself and continued on with the lesson."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    Write-Host "WARNING: target sentence was not found; document left unchanged."
} else {
    Write-Host "Merged the split runs back into a single run."
}
